$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (estoque_atualizado: -252 -> -254)
$ws.Range("G2").Value = -254

# Row 3 (estoque_atualizado: -109 -> -110)
$ws.Range("G3").Value = -110

# Row 5 (estoque_atualizado: -74 -> -75, media_vendas: 1.06 -> 1.05)
$ws.Range("G5").Value = -75
$ws.Range("H5").Value = 1.05

# Row 7 (estoque_atualizado: -61 -> -63, desvio_padrao: 0.17 -> 0.16)
$ws.Range("G7").Value = -63
$ws.Range("I7").Value = 0.16

# Row 9 (estoque_atualizado: -1251 -> -1252, media_vendas: 1.06 -> 1.05, desvio_padrao: 0.31 -> 0.26)
$ws.Range("G9").Value = -1252
$ws.Range("H9").Value = 1.05
$ws.Range("I9").Value = 0.26
